# Adds a new "2022-Q3" quarterly sheet (right after "总计") and records its
# summary line in the "总计" sheet, matching commit "feat: add 2022-Q3 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q1"), so every other quarterly sheet
#    keeps its existing name -> data mapping untouched.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Header row (same headers/order as the other quarterly sheets).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Fund rows for 2022-Q3 (A=index, B=code, C=name, D=scale, E=equity position,
# F=weight, G=market value held (CNY100M), H=weight rank). D..G mirror the
# source workbook where these numeric-looking figures are stored as text.
$rows = @(
    @("010041", "嘉实港股优势混合A", "40.43", "89.69", "3.22", "1.3018", 9),
    @("009983", "永赢港股通品质生活慧选混合", "9.17", "60.68", "3.97", "0.3640", 5),
    @("010042", "嘉实港股优势混合C", "5.20", "89.69", "3.22", "0.1674", 9),
    @("011315", "永赢港股通优质成长一年混合", "3.42", "60.58", "3.77", "0.1289", 8),
    @("004317", "前海开源沪港深裕鑫灵活配置混合C", "1.79", "70.17", "3.10", "0.0555", 1),
    @("004316", "前海开源沪港深裕鑫灵活配置混合A", "1.77", "70.17", "3.10", "0.0549", 1),
    @("161124", "易方达香港恒生综合小型股指数（QDII-LOF）A", "0.20", "91.61", "2.28", "0.0046", 1),
    @("006263", "易方达香港恒生综合小型股指数（QDII-LOF）C", "0.05", "91.61", "2.28", "0.0011", 1)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $excelRow = $r + 2
    $data = $rows[$r]

    $idxCell = $newSheet.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $codeCell = $newSheet.Cells.Item($excelRow, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $data[0]

    $newSheet.Cells.Item($excelRow, 3).Value = $data[1]

    $scaleCell = $newSheet.Cells.Item($excelRow, 4)
    $scaleCell.NumberFormat = "@"
    $scaleCell.Value = $data[2]

    $posCell = $newSheet.Cells.Item($excelRow, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $data[3]

    $weightCell = $newSheet.Cells.Item($excelRow, 6)
    $weightCell.NumberFormat = "@"
    $weightCell.Value = $data[4]

    $mvCell = $newSheet.Cells.Item($excelRow, 7)
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $data[5]

    $newSheet.Cells.Item($excelRow, 8).Value = $data[6]
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Insert a new summary row into "总计" for the 2022-Q3 quarter, right
#    above the existing "2022-Q1" row (pushing the rest down by one).
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).ClearFormats()

$tIdx = $totalSheet.Cells.Item(2, 1)
$tIdx.Value = 0
$tIdx.Font.Bold = $true
$tIdx.HorizontalAlignment = -4108
$tIdx.VerticalAlignment = -4160
$tIdx.Borders.LineStyle = 1

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 2.08

$totalSheet.Range("A1").Select()
